$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 stays "Yard" (re-pointed to a different shared-string slot upstream,
# value itself is unchanged)
$ws.Range("B2").Value = "Yard"

# Column P ("Beacon Info") gets populated with station/track beacon codes
$ws.Range("P3").Value = "B1,2,Pioneer"
$ws.Range("P5").Value = "B2,2,Pioneer"
$ws.Range("P10").Value = "B3,2,Edgebrook"
$ws.Range("P12").Value = "B4,2,Edgebrook"
$ws.Range("P17").Value = "B5,3,Lebron"
$ws.Range("P19").Value = "B6,3,Lebron"
$ws.Range("P23").Value = "B7,3,Whited"
$ws.Range("P25").Value = "B8,3,Whited"
$ws.Range("P32").Value = "B9,2,South Bank"
$ws.Range("P34").Value = "B10,2,South Bank"
$ws.Range("P38").Value = "T1"
$ws.Range("P40").Value = "B11,1,Central"
$ws.Range("P42").Value = "B12,1,Central"
$ws.Range("P49").Value = "B13,1,Inglewood"
$ws.Range("P51").Value = "B14,1,Inglewood"
$ws.Range("P58").Value = "B15,1,Overbrook"
$ws.Range("P59").Value = "T2"
$ws.Range("P60").Value = "B16,1,Overbrook"
$ws.Range("P66").Value = "B17,1,Glenbury"
$ws.Range("P68").Value = "B18,1,Glenbury"
$ws.Range("P74").Value = "B19,1,Dormont"
$ws.Range("P76").Value = "B20,1,Dormont"
$ws.Range("P78").Value = "B21,3,Mt Lebanon"
$ws.Range("P80").Value = "B22,3,Mt Lebanon"
$ws.Range("P89").Value = "B23,2,Mt Poplar"
$ws.Range("P91").Value = "B24,2,Mt Poplar"
$ws.Range("P97").Value = "B25,2,Castle Shannon"
$ws.Range("P99").Value = "B26,2,Castle Shannon"
$ws.Range("P106").Value = "B27,1,Dormont"
$ws.Range("P108").Value = "B28,1,Dormont"
$ws.Range("P115").Value = "B29,1,Glenbury"
$ws.Range("P117").Value = "B30,1,Glenbury"
$ws.Range("P124").Value = "B31,1,Overbrook"
$ws.Range("P126").Value = "B32,1,Overbrook"
$ws.Range("P133").Value = "B33,2,Inglewood"
$ws.Range("P135").Value = "B34,2,Inglewood"
$ws.Range("P142").Value = "B35,1,Central"
$ws.Range("P144").Value = "B36,1,Central"
$ws.Range("P145").Value = "T3"

# View-state tweaks: scroll position + selection (best effort - the
# window's horizontal scroll anchor isn't separately exposed beyond the
# active cell in this host, so we drive it via selection)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 8
$ws.Range("P1:P152").Select()
